$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp
$ws.Range("A1").Value = "Datos actualizados a 25 de Julio de 2020 a las 12:51"

# --- Row 4: Estados Unidos (simple data refresh) ---
$ws.Range("B4").Value = 4248759
$ws.Range("C4").Value = 432
$ws.Range("E4").Value = 2071900
$ws.Range("G4").Value = 8
$ws.Range("H4").Value = 148498

# --- Row 14: Iran (simple data refresh) ---
$ws.Range("B14").Value = 288839
$ws.Range("C14").Value = 2316
$ws.Range("D14").Value = 251319
$ws.Range("E14").Value = 22036
$ws.Range("G14").Value = 195
$ws.Range("H14").Value = 15484

# --- Row 25: Catar (simple data refresh) ---
$ws.Range("B25").Value = 109036
$ws.Range("C25").Value = 398
$ws.Range("D25").Value = 105750
$ws.Range("E25").Value = 3122

# --- Rows 47-49: Rumania moves up above Guatemala and Polonia ---
$ws.Range("A47").Value = "Rumania"
$ws.Range("B47").Value = 43678
$ws.Range("C47").Value = 1284
$ws.Range("D47").Value = 25373
$ws.Range("E47").Value = 16140
$ws.Range("G47").Value = 15
$ws.Range("H47").Value = 2165

$ws.Range("A48").Value = "Guatemala"
$ws.Range("B48").Value = 43283
$ws.Range("C48").Value = 0
$ws.Range("D48").Value = 30150
$ws.Range("E48").Value = 11464
$ws.Range("G48").Value = 0
$ws.Range("H48").Value = 1669

$ws.Range("A49").Value = "Polonia"
$ws.Range("B49").Value = 42622
$ws.Range("C49").Value = 584
$ws.Range("D49").Value = 32419
$ws.Range("E49").Value = 8539
$ws.Range("G49").Value = 9
$ws.Range("H49").Value = 1664

# --- Rows 73-75: El Salvador moves up above Corea del Sur and Australia ---
$ws.Range("A73").Value = "El Salvador"
$ws.Range("B73").Value = 14221
$ws.Range("C73").Value = 429
$ws.Range("D73").Value = 7549
$ws.Range("E73").Value = 6282
$ws.Range("G73").Value = 11
$ws.Range("H73").Value = 390

$ws.Range("A74").Value = "Corea del Sur"
$ws.Range("B74").Value = 14092
$ws.Range("C74").Value = 113
$ws.Range("D74").Value = 12866
$ws.Range("E74").Value = 928
$ws.Range("G74").Value = 0
$ws.Range("H74").Value = 298

$ws.Range("A75").Value = "Australia"
$ws.Range("B75").Value = 13948
$ws.Range("C75").Value = 353
$ws.Range("D75").Value = 8929
$ws.Range("E75").Value = 4874
$ws.Range("G75").Value = 6
$ws.Range("H75").Value = 145

# --- Row 87: Consejo Danes para los Refugiados (simple data refresh) ---
$ws.Range("B87").Value = 8801
$ws.Range("C87").Value = 34
$ws.Range("D87").Value = 5305
$ws.Range("E87").Value = 3292
$ws.Range("G87").Value = 3
$ws.Range("H87").Value = 204

# --- Rows 210-211: Groenlandia moves above Islas Malvinas (values identical, only names swap) ---
$ws.Range("A210").Value = "Groenlandia"
$ws.Range("A211").Value = "Islas Malvinas"

# --- Row 104: Zambia (simple data refresh) ---
$ws.Range("E104").Value = 2040
$ws.Range("G104").Value = 3
$ws.Range("H104").Value = 139

# --- Row 125: Eslovenia (simple data refresh) ---
$ws.Range("B125").Value = 2066
$ws.Range("C125").Value = 14
$ws.Range("E125").Value = 272
$ws.Range("G125").Value = 1
$ws.Range("H125").Value = 116

# --- Row 145: Uganda (simple data refresh) ---
$ws.Range("B145").Value = 1103
$ws.Range("C145").Value = 14
$ws.Range("E145").Value = 127
